$p = $ppt.ActivePresentation
$s = $p.Slides.Item(36)
$shp = $s.Shapes.Item("Group 11")
$shp.Top = 290.368346456693
